$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.520.06"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.752.74"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.03"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.88"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.751.51"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.45"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("E13").Value = "  -7.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.09"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.384.20"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.752.89"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.535.52"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.99"
$ws.Range("E18").Value = "  -4.84%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.75"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "464.29"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("E23").Value = "  -3.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.04"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  -3.59%  "
$ws.Range("E26").Value = "  -3.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.93"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  -4.89%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.899.80"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("E31").Value = "  -5.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.30"
$ws.Range("E32").Value = "  -3.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.96"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.16"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.704.68"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.100"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("E39").Value = "  -9.56%  "
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.00"
$ws.Range("E45").Value = "  +9.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.302"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.91"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.56"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.73"
$ws.Range("E50").Value = "  +3.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "387.15"
